# Optuna Attempt (go back with original)
#
# Updates forecast metrics on the "Forecast Comparison" sheet (Seasonality
# Index / Inventory Coverage / MyForecast values) and the derived rollup
# totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# MyForecast (column D)
$wsForecast.Range("D3").Value = 119

# Inventory Coverage (column H)
$wsForecast.Range("H3").Value = 5.81
$wsForecast.Range("H4").Value = 5.68
$wsForecast.Range("H5").Value = 4.56
$wsForecast.Range("H6").Value = 3.52
$wsForecast.Range("H7").Value = 2.55
$wsForecast.Range("H8").Value = 1.51
$wsForecast.Range("H9").Value = 0.5

# Seasonality Index (column L)
$wsForecast.Range("L2").Value = 1.1
$wsForecast.Range("L3").Value = 1.04
$wsForecast.Range("L4").Value = 0.99
$wsForecast.Range("L5").Value = 0.82
$wsForecast.Range("L6").Value = 0.97
$wsForecast.Range("L7").Value = 0.99
$wsForecast.Range("L8").Value = 0.99
$wsForecast.Range("L9").Value = 1.06
$wsForecast.Range("L10").Value = 0.86
$wsForecast.Range("L11").Value = 1.15
$wsForecast.Range("L12").Value = 0.86
$wsForecast.Range("L13").Value = 1.05
$wsForecast.Range("L14").Value = 0.91
$wsForecast.Range("L15").Value = 0.96
$wsForecast.Range("L16").Value = 1.02
$wsForecast.Range("L17").Value = 0.8

# --- Summary sheet ---------------------------------------------------------
# These cells hold their numbers as text (e.g. "1656"), so force a text
# number format before writing the replacement digits - otherwise Excel's
# COM layer infers a numeric-looking string and stores it as a Number.

$summaryCells = @(
    @{ Addr = "B9";  Val = "1661" },
    @{ Addr = "B10"; Val = "852"  },
    @{ Addr = "B11"; Val = "428"  },
    @{ Addr = "B12"; Val = "120"  }
)

foreach ($item in $summaryCells) {
    $rng = $wsSummary.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
}
